# Auto-generated: updates Typhon Profits market-price columns (H:N) per the
# scheduled runner's latest Universalis snapshot. Cells whose computed
# Leve-profit (M) or (N) value is now blank/N-A are cleared; cells that newly
# resolve to a number are created.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 361.07144
$ws.Range("I4").Value = 150.55556
$ws.Range("J4").Value = 740
$ws.Range("K4").Value = 150.55556
$ws.Range("L4").Value = 740
$ws.Range("M4").Value = -36.55556000000001
$ws.Range("N4").Value = -968
$ws.Range("H19").Value = 2564.3635
$ws.Range("I19").Value = 4557.6
$ws.Range("K19").Value = 4557.6
$ws.Range("M19").Value = -4382.6
$ws.Range("H32").Value = 766.6667
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 840
$ws.Range("K32").Value = 400
$ws.Range("L32").Value = 840
$ws.Range("M32").Value = -74
$ws.Range("N32").Value = -1492
$ws.Range("H70").Value = 1146.3334
$ws.Range("I70").Value = 1102.125
$ws.Range("K70").Value = 3306.375
$ws.Range("M70").Value = -3036.375
$ws.Range("H73").Value = 1146.3334
$ws.Range("I73").Value = 1102.125
$ws.Range("K73").Value = 3306.375
$ws.Range("M73").Value = -2370.375
$ws.Range("H76").Value = 4632554.5
$ws.Range("I76").Value = 3160
$ws.Range("K76").Value = 3160
$ws.Range("M76").Value = -2845
$ws.Range("H79").Value = 4632554.5
$ws.Range("I79").Value = 3160
$ws.Range("K79").Value = 3160
$ws.Range("M79").Value = -2068
$ws.Range("H106").Value = 1874.9445
$ws.Range("I106").Value = 1621.1875
$ws.Range("K106").Value = 1621.1875
$ws.Range("M106").Value = -990.1875
$ws.Range("H129").Value = 837.7778
$ws.Range("J129").Value = 850
$ws.Range("L129").Value = 2550
$ws.Range("N129").Value = -12550
$ws.Range("H132").Value = 3034.6428
$ws.Range("I132").Value = 3058.8
$ws.Range("J132").Value = 2833.3333
$ws.Range("K132").Value = 9176.400000000001
$ws.Range("L132").Value = 8499.999899999999
$ws.Range("M132").Value = -6646.400000000001
$ws.Range("N132").Value = -13559.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23091.021
$ws.Range("I32").Value = 25406.627
$ws.Range("K32").Value = 25406.627
$ws.Range("M32").Value = -25119.627
$ws.Range("H61").Value = 4032.6
$ws.Range("I61").Value = 2898.7778
$ws.Range("K61").Value = 2898.7778
$ws.Range("M61").Value = -2686.7778
$ws.Range("H74").Value = 2023.0435
$ws.Range("I74").Value = 2190.5557
$ws.Range("K74").Value = 2190.5557
$ws.Range("M74").Value = -1316.5557
$ws.Range("H77").Value = 2023.0435
$ws.Range("I77").Value = 2190.5557
$ws.Range("K77").Value = 10952.7785
$ws.Range("M77").Value = -6584.7785
$ws.Range("H102").Value = 2526.2104
$ws.Range("I102").Value = 978.38464
$ws.Range("J102").Value = 5879.8335
$ws.Range("K102").Value = 978.38464
$ws.Range("L102").Value = 5879.8335
$ws.Range("M102").Value = 643.61536
$ws.Range("N102").Value = -9123.833500000001
$ws.Range("H122").Value = 1449.3572
$ws.Range("I122").Value = 1533.6666
$ws.Range("K122").Value = 4600.9998
$ws.Range("M122").Value = -2150.9998
$ws.Range("H132").Value = 22922.6
$ws.Range("I132").Value = 2782.6365
$ws.Range("J132").Value = 38746.855
$ws.Range("K132").Value = 8347.9095
$ws.Range("L132").Value = 116240.565
$ws.Range("M132").Value = -5817.9095
$ws.Range("N132").Value = -121300.565
$ws.Range("H136").Value = 4032.6
$ws.Range("I136").Value = 2898.7778
$ws.Range("K136").Value = 8696.3334
$ws.Range("M136").Value = -6146.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1050.6666
$ws.Range("I20").Value = 1125.5
$ws.Range("J20").Value = 901
$ws.Range("K20").Value = 1125.5
$ws.Range("L20").Value = 901
$ws.Range("M20").Value = -878.5
$ws.Range("N20").Value = -1395
$ws.Range("H107").Value = 2584.8
$ws.Range("I107").Value = 1477.75
$ws.Range("K107").Value = 1477.75
$ws.Range("M107").Value = 442.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 10000
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 109.07692
$ws.Range("I6").Value = 83.27273
$ws.Range("J6").Value = 251
$ws.Range("K6").Value = 249.81819
$ws.Range("L6").Value = 753
$ws.Range("M6").Value = -136.81819
$ws.Range("N6").Value = -979
$ws.Range("H107").Value = 16979
$ws.Range("J107").Value = 374.8
$ws.Range("L107").Value = 1124.4
$ws.Range("N107").Value = -4964.4
$ws.Range("H131").Value = 760.91
$ws.Range("I131").Value = 367.5
$ws.Range("J131").Value = 777.30206
$ws.Range("K131").Value = 1102.5
$ws.Range("L131").Value = 2331.90618
$ws.Range("M131").Value = 3937.5
$ws.Range("N131").Value = -12411.90618
$ws.Range("H132").Value = 874.6667
$ws.Range("I132").Value = 874.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7872.0003
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -5342.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3376.04
$ws.Range("I80").Value = 3094.4546
$ws.Range("J80").Value = 3597.2856
$ws.Range("K80").Value = 3094.4546
$ws.Range("L80").Value = 3597.2856
$ws.Range("M80").Value = -2096.4546
$ws.Range("N80").Value = -5593.2856
$ws.Range("H83").Value = 3376.04
$ws.Range("I83").Value = 3094.4546
$ws.Range("J83").Value = 3597.2856
$ws.Range("K83").Value = 15472.273
$ws.Range("L83").Value = 17986.428
$ws.Range("M83").Value = -10480.273
$ws.Range("N83").Value = -27970.428
$ws.Range("H122").Value = 1528
$ws.Range("I122").Value = 1491.5714
$ws.Range("K122").Value = 4474.7142
$ws.Range("M122").Value = -2024.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1866.75
$ws.Range("I93").Value = 1851.9445
$ws.Range("K93").Value = 1851.9445
$ws.Range("M93").Value = -603.9445000000001
$ws.Range("H122").Value = 1228915.5
$ws.Range("J122").Value = 4077.7778
$ws.Range("L122").Value = 12233.3334
$ws.Range("N122").Value = -17133.3334
$ws.Range("H132").Value = 1940.55
$ws.Range("I132").Value = 1327.5385
$ws.Range("J132").Value = 3079
$ws.Range("K132").Value = 3982.6155
$ws.Range("L132").Value = 9237
$ws.Range("M132").Value = -1452.6155
$ws.Range("N132").Value = -14297

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = $null
$ws.Range("H136").Value = 29413112
$ws.Range("I136").Value = 41667572
$ws.Range("J136").Value = 2410.3
$ws.Range("K136").Value = 125002716
$ws.Range("L136").Value = 7230.900000000001
$ws.Range("M136").Value = -125000166
$ws.Range("N136").Value = -12330.9

